$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("2025-05-10", "38", "37.05", "0.98", "0.265", "0.09", "5,311", "7,951", "8,001", "7.2617")

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(70, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
}
